# Feedback form basic layout
# Adds a new "get user feedback" (WS-FED-03) service row to the ServicesList sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServicesList")

$row = 40
$prev = 39

# Column F needs to hold the literal text "false" (matching the rest of the
# column, which is text, not a real boolean) -- enter it with a leading
# apostrophe so it isn't auto-converted to a Boolean.
$ws.Cells.Item($row, 6).Value = "'false"

# Copy the formatting (styles/borders) of the row above into the new row,
# matching the existing table's look (columns B:L).
$ws.Range("B$prev`:L$prev").Copy() | Out-Null
$ws.Range("B$row`:L$row").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column B: feedback (group/category)
$ws.Cells.Item($row, 2).Value = "feedback"
# Column C: Service Name
$ws.Cells.Item($row, 3).Value = "get user feedback"
# Column D: ServiceCode
$ws.Cells.Item($row, 4).Value = "WS-FED-03"
# Column E: queryId
$ws.Cells.Item($row, 5).Value = "app.feedback.get"
# Column G: BasePath
$ws.Cells.Item($row, 7).Value = "feedback"
# Column H: servicePath
$ws.Cells.Item($row, 8).Value = "/get"
# Column I: ServiceType
$ws.Cells.Item($row, 9).Value = "POST"
# Columns J, K, L left blank (Priority, Coding, Testing)

# Column M: generated SQL insert statement
$ws.Cells.Item($row, 13).Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D40,"'',''CONNON_CONFIG'', 0, ''",C40,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'

# Column N: generated mapping annotation
$ws.Cells.Item($row, 14).Formula = '=_xlfn.CONCAT(IF(I40="GET","@GetMapping(",IF(I40="POST","@PostMapping(",IF(I40="DELETE","@DeleteMapping(",IF(I40="PUT","@PutMapping(","")))),CHAR(34),H40,CHAR(34),")")'

# Column O: generated ServiceInfo annotation
$ws.Cells.Item($row, 15).Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D40,,CHAR(34),", serviceName = ",CHAR(34),C40,CHAR(34), ", queryId = ",CHAR(34),E40,CHAR(34),", logActivity =",F40,")")'

# Update selection to match the author's final cursor position
$ws.Range("N40:O40").Select()
